$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 761.75
$ws.Range("I111").Value = 804.8333
$ws.Range("J111").Value = 632.5
$ws.Range("K111").Value = 2414.4999
$ws.Range("L111").Value = 1897.5
$ws.Range("M111").Value = 652.5001000000002
$ws.Range("N111").Value = -8031.5

$ws.Range("H113").Value = 1523.4615
$ws.Range("J113").Value = 1555.091
$ws.Range("L113").Value = 1555.091
$ws.Range("N113").Value = -8063.091

$ws.Range("H125").Value = 2775
$ws.Range("I125").Value = 2888
$ws.Range("J125").Value = 2436
$ws.Range("K125").Value = 25992
$ws.Range("L125").Value = 21924
$ws.Range("M125").Value = -23532
$ws.Range("N125").Value = -26844

$ws.Range("H127").Value = 741.6667
$ws.Range("I127").Value = 470
$ws.Range("J127").Value = 1149.1666
$ws.Range("K127").Value = 1410
$ws.Range("L127").Value = 3447.4998
$ws.Range("M127").Value = 3550
$ws.Range("N127").Value = -13367.4998

$ws.Range("H132").Value = 190871.45
$ws.Range("I132").Value = 215110.25
$ws.Range("J132").Value = 1000.8333
$ws.Range("K132").Value = 645330.75
$ws.Range("L132").Value = 3002.4999
$ws.Range("M132").Value = -642800.75
$ws.Range("N132").Value = -8062.4999

$ws.Range("H138").Value = 4297.75
$ws.Range("I138").Value = 5090.9688
$ws.Range("J138").Value = 3504.5312
$ws.Range("K138").Value = 15272.9064
$ws.Range("L138").Value = 10513.5936
$ws.Range("M138").Value = -10132.9064
$ws.Range("N138").Value = -20793.5936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7075.523
$ws.Range("I32").Value = 7481.324
$ws.Range("K32").Value = 7481.324
$ws.Range("M32").Value = -7194.324

$ws.Range("H45").Value = 887.4286
$ws.Range("I45").Value = 887.4286
$ws.Range("K45").Value = 887.4286
$ws.Range("M45").Value = -510.4286

$ws.Range("H61").Value = 17546124
$ws.Range("I61").Value = 22224436
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 22224436
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -22224224
$ws.Range("N61").Value = -2874

$ws.Range("H132").Value = 3497.1765
$ws.Range("I132").Value = 3976.889
$ws.Range("K132").Value = 11930.667
$ws.Range("M132").Value = -9400.667000000001

$ws.Range("H136").Value = 17546124
$ws.Range("I136").Value = 22224436
$ws.Range("J136").Value = 2450
$ws.Range("K136").Value = 66673308
$ws.Range("L136").Value = 7350
$ws.Range("M136").Value = -66670758
$ws.Range("N136").Value = -12450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2414.4614
$ws.Range("I20").Value = 2587.5557
$ws.Range("J20").Value = 2025
$ws.Range("K20").Value = 2587.5557
$ws.Range("L20").Value = 2025
$ws.Range("M20").Value = -2340.5557
$ws.Range("N20").Value = -2519

$ws.Range("H86").Value = 2551.7646
$ws.Range("I86").Value = 2054.5454
$ws.Range("J86").Value = 3463.3333
$ws.Range("K86").Value = 2054.5454
$ws.Range("L86").Value = 3463.3333
$ws.Range("M86").Value = -931.5454
$ws.Range("N86").Value = -5709.3333

$ws.Range("H89").Value = 2551.7646
$ws.Range("I89").Value = 2054.5454
$ws.Range("J89").Value = 3463.3333
$ws.Range("K89").Value = 10272.727
$ws.Range("L89").Value = 17316.6665
$ws.Range("M89").Value = -4656.726999999999
$ws.Range("N89").Value = -28548.6665

$ws.Range("H128").Value = 2996.6667
$ws.Range("I128").Value = 2996.6667
$ws.Range("K128").Value = 8990.000100000001
$ws.Range("M128").Value = -6500.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2601.25
$ws.Range("I62").Value = 2502.5
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 2502.5
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -1878.5
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 2601.25
$ws.Range("I65").Value = 2502.5
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 12512.5
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -9392.5
$ws.Range("N65").Value = -19740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1500
$ws.Range("J9").Value = 1500
$ws.Range("L9").Value = 4500
$ws.Range("N9").Value = -4948

$ws.Range("H34").Value = 90909976
$ws.Range("J34").Value = 100000960
$ws.Range("L34").Value = 300002880
$ws.Range("N34").Value = -300003048

$ws.Range("H113").Value = 592.975
$ws.Range("I113").Value = 513.86365
$ws.Range("J113").Value = 689.6667
$ws.Range("K113").Value = 1541.59095
$ws.Range("L113").Value = 2069.0001
$ws.Range("M113").Value = 628.40905
$ws.Range("N113").Value = -6409.0001

$ws.Range("H114").Value = 1579.3334
$ws.Range("I114").Value = 282.18182
$ws.Range("J114").Value = 2676.923
$ws.Range("K114").Value = 846.54546
$ws.Range("L114").Value = 8030.768999999999
$ws.Range("M114").Value = 2407.45454
$ws.Range("N114").Value = -14538.769

$ws.Range("H131").Value = 3217.3333
$ws.Range("J131").Value = 2568.6562
$ws.Range("L131").Value = 7705.9686
$ws.Range("N131").Value = -17785.9686

$ws.Range("H134").Value = 4954.65
$ws.Range("J134").Value = 6799.25
$ws.Range("L134").Value = 20397.75
$ws.Range("N134").Value = -30537.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 551.73334
$ws.Range("I107").Value = 278.22223
$ws.Range("K107").Value = 278.22223
$ws.Range("M107").Value = 1641.77777

$ws.Range("H122").Value = 2386.0356
$ws.Range("I122").Value = 1585.762
$ws.Range("J122").Value = 4786.857
$ws.Range("K122").Value = 4757.286
$ws.Range("L122").Value = 14360.571
$ws.Range("M122").Value = -2307.286
$ws.Range("N122").Value = -19260.571

$ws.Range("H132").Value = 1605.28
$ws.Range("I132").Value = 1320.25
$ws.Range("J132").Value = 2745.4
$ws.Range("K132").Value = 3960.75
$ws.Range("L132").Value = 8236.200000000001
$ws.Range("M132").Value = -1430.75
$ws.Range("N132").Value = -13296.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6261310.5
$ws.Range("I96").Value = 12501851
$ws.Range("J96").Value = 20770.5
$ws.Range("K96").Value = 12501851
$ws.Range("L96").Value = 20770.5
$ws.Range("M96").Value = -12500478
$ws.Range("N96").Value = -23516.5

$ws.Range("H113").Value = 622.6667
$ws.Range("I113").Value = 624
$ws.Range("J113").Value = 621.7143
$ws.Range("K113").Value = 1872
$ws.Range("L113").Value = 1865.1429
$ws.Range("M113").Value = 298
$ws.Range("N113").Value = -6205.1429

$ws.Range("H126").Value = 2819.9333
$ws.Range("I126").Value = 4112.222
$ws.Range("J126").Value = 881.5
$ws.Range("K126").Value = 12336.666
$ws.Range("L126").Value = 2644.5
$ws.Range("M126").Value = -9866.665999999999
$ws.Range("N126").Value = -7584.5

$ws.Range("H132").Value = 2538.8096
$ws.Range("I132").Value = 2575.875
$ws.Range("J132").Value = 2420.2
$ws.Range("K132").Value = 7727.625
$ws.Range("L132").Value = 7260.599999999999
$ws.Range("M132").Value = -5197.625
$ws.Range("N132").Value = -12320.6
